$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CasesTab" query in B2 no longer returns a Cohort column - the trailing
# "OPTIONAL MATCH (co:cohort)..." / "coalesce(co.cohort_description, '') AS `Cohort`"
# pairing is dropped, leaving `Response to Treatment` as the final RETURN column.
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Belgian Malinois', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder, Prostate', 'Bladder, Urethra', 'Bladder, Urethra, Prostate', 'Urethra, Prostate'] and diag.best_response in ['Not Determined', 'Partial Response']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Restore the top-left/selection view to its post-edit state (Excel re-homes the
# window/selection to the top of the sheet after this timing-issue fix).
$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
